# aggiornamento fino a 27/05
# Append daily rows (date, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
# for 14 new days (14/05/2021 - 27/05/2021) to the bottom of the existing data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 255
$newData = @(
    @(44330, 7, 28, 84.89221720280145),
    @(44331, 2, 20, 60.63729800200103),
    @(44332, 5, 19, 57.60543310190099),
    @(44333, 3, 20, 60.63729800200103),
    @(44334, 1, 20, 60.63729800200103),
    @(44335, 0, 19, 57.60543310190099),
    @(44336, 0, 18, 54.57356820180092),
    @(44337, 6, 17, 51.54170330170087),
    @(44338, 4, 19, 57.60543310190099),
    @(44339, 0, 14, 42.44610860140072),
    @(44340, 4, 15, 45.47797350150077),
    @(44341, 4, 18, 54.57356820180092),
    @(44342, 1, 19, 57.60543310190099),
    @(44343, 4, 23, 69.73289270230119)
)

$r = $lastRow
foreach ($entry in $newData) {
    $r = $r + 1

    # Copy the date cell above so the new date cell inherits the same
    # style (border/bold/center/top alignment + date number format).
    $ws.Range("A$lastRow").Copy($ws.Range("A$r"))

    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}
